$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2815276.8
$ws.Range("J17").Value = 2865540.8
$ws.Range("L17").Value = 8596622.399999999
$ws.Range("N17").Value = -8596958.399999999

$ws.Range("H33").Value = 3367989.5
$ws.Range("I33").Value = 1164.7368
$ws.Range("K33").Value = 1164.7368
$ws.Range("M33").Value = -935.7367999999999

$ws.Range("H107").Value = 12500472
$ws.Range("I107").Value = 13889202
$ws.Range("K107").Value = 13889202
$ws.Range("M107").Value = -13887282

$ws.Range("H129").Value = 978.4776000000001
$ws.Range("J129").Value = 1065.7455
$ws.Range("L129").Value = 3197.2365
$ws.Range("N129").Value = -13197.2365

$ws.Range("H132").Value = 1699.5769
$ws.Range("I132").Value = 1745.174
$ws.Range("J132").Value = 1350
$ws.Range("K132").Value = 5235.522
$ws.Range("L132").Value = 4050
$ws.Range("M132").Value = -2705.522
$ws.Range("N132").Value = -9110

$ws.Range("H137").Value = 1416.2354
$ws.Range("I137").Value = 1412.8334
$ws.Range("J137").Value = 1424.4
$ws.Range("K137").Value = 4238.5002
$ws.Range("L137").Value = 4273.200000000001
$ws.Range("M137").Value = -1688.5002
$ws.Range("N137").Value = -9373.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2850474
$ws.Range("I102").Value = 3368299
$ws.Range("K102").Value = 3368299
$ws.Range("M102").Value = -3366677

$ws.Range("H122").Value = 1351186
$ws.Range("I122").Value = 1604142.8
$ws.Range("J122").Value = 2083.3333
$ws.Range("K122").Value = 4812428.4
$ws.Range("L122").Value = 6249.999899999999
$ws.Range("M122").Value = -4809978.4
$ws.Range("N122").Value = -11149.9999

$ws.Range("H123").Value = 29220
$ws.Range("J123").Value = 29220
$ws.Range("L123").Value = 29220
$ws.Range("N123").Value = -39020

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11696.458
$ws.Range("I20").Value = 1170.6
$ws.Range("J20").Value = 19214.928
$ws.Range("K20").Value = 1170.6
$ws.Range("L20").Value = 19214.928
$ws.Range("M20").Value = -923.5999999999999
$ws.Range("N20").Value = -19708.928

$ws.Range("H94").Value = 1650.5294
$ws.Range("I94").Value = 914.4545000000001
$ws.Range("K94").Value = 914.4545000000001
$ws.Range("M94").Value = -463.4545000000001

$ws.Range("H132").Value = 40032.25
$ws.Range("J132").Value = 40032.25
$ws.Range("L132").Value = 40032.25
$ws.Range("N132").Value = -50152.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1542.5312
$ws.Range("I58").Value = 974.8125
$ws.Range("J58").Value = 2110.25
$ws.Range("K58").Value = 974.8125
$ws.Range("L58").Value = 2110.25
$ws.Range("M58").Value = -771.8125
$ws.Range("N58").Value = -2516.25

$ws.Range("H134").Value = 3130.6765
$ws.Range("I134").Value = 3378.8076
$ws.Range("J134").Value = 2324.25
$ws.Range("K134").Value = 10136.4228
$ws.Range("L134").Value = 6972.75
$ws.Range("M134").Value = -7601.4228
$ws.Range("N134").Value = -12042.75

$ws.Range("H136").Value = 1542.5312
$ws.Range("I136").Value = 974.8125
$ws.Range("J136").Value = 2110.25
$ws.Range("K136").Value = 2924.4375
$ws.Range("L136").Value = 6330.75
$ws.Range("M136").Value = -374.4375
$ws.Range("N136").Value = -11430.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 436.9
$ws.Range("I60").Value = 254.14285
$ws.Range("J60").Value = 863.3333
$ws.Range("K60").Value = 762.4285500000001
$ws.Range("L60").Value = 2589.9999
$ws.Range("M60").Value = -511.4285500000001
$ws.Range("N60").Value = -3091.9999

$ws.Range("H69").Value = 1382.1538
$ws.Range("J69").Value = 1718.7778
$ws.Range("L69").Value = 5156.3334
$ws.Range("N69").Value = -6778.3334

$ws.Range("H72").Value = 1382.1538
$ws.Range("J72").Value = 1718.7778
$ws.Range("L72").Value = 15469.0002
$ws.Range("N72").Value = -23581.0002

$ws.Range("H113").Value = 2222908.5
$ws.Range("I113").Value = 2632277.8
$ws.Range("K113").Value = 7896833.399999999
$ws.Range("M113").Value = -7894663.399999999

$ws.Range("H122").Value = 2329.2744
$ws.Range("I122").Value = 660
$ws.Range("J122").Value = 2510.7173
$ws.Range("K122").Value = 5940
$ws.Range("L122").Value = 22596.4557
$ws.Range("M122").Value = -3490
$ws.Range("N122").Value = -27496.4557

$ws.Range("H129").Value = 1185.9231
$ws.Range("I129").Value = 991.7
$ws.Range("J129").Value = 1833.3334
$ws.Range("K129").Value = 2975.1
$ws.Range("L129").Value = 5500.0002
$ws.Range("M129").Value = 2024.9
$ws.Range("N129").Value = -15500.0002

$ws.Range("H137").Value = 10177.553
$ws.Range("I137").Value = 8201.1875
$ws.Range("J137").Value = 11197.613
$ws.Range("K137").Value = 24603.5625
$ws.Range("L137").Value = 33592.839
$ws.Range("M137").Value = -19503.5625
$ws.Range("N137").Value = -43792.839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 1099.8334
$ws.Range("I97").Value = 1059.8
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 1059.8
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -563.8
$ws.Range("N97").Value = -2292

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2495331.5
$ws.Range("I122").Value = 3602763
$ws.Range("J122").Value = 3611
$ws.Range("K122").Value = 10808289
$ws.Range("L122").Value = 10833
$ws.Range("M122").Value = -10805839
$ws.Range("N122").Value = -15733

$ws.Range("H123").Value = 26925.834
$ws.Range("J123").Value = 26925.834
$ws.Range("L123").Value = 26925.834
$ws.Range("N123").Value = -31825.834

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 603.86957
$ws.Range("I46").Value = 415.9
$ws.Range("J46").Value = 748.46155
$ws.Range("K46").Value = 415.9
$ws.Range("L46").Value = 748.46155
$ws.Range("M46").Value = -227.9
$ws.Range("N46").Value = -1124.46155

$ws.Range("H61").Value = 1589.4445
$ws.Range("I61").Value = 1593.5714
$ws.Range("K61").Value = 1593.5714
$ws.Range("M61").Value = -1391.5714

$ws.Range("H93").Value = 1224.75
$ws.Range("I93").Value = 1199.75
$ws.Range("J93").Value = 1249.75
$ws.Range("K93").Value = 1199.75
$ws.Range("L93").Value = 1249.75
$ws.Range("M93").Value = 48.25
$ws.Range("N93").Value = -3745.75

$ws.Range("H113").Value = 1589.4445
$ws.Range("I113").Value = 1593.5714
$ws.Range("K113").Value = 1593.5714
$ws.Range("M113").Value = 576.4286

$ws.Range("H122").Value = 3395240
$ws.Range("I122").Value = 4466453.5
$ws.Range("K122").Value = 13399360.5
$ws.Range("M122").Value = -13396910.5

$ws.Range("H136").Value = 5072.7383
$ws.Range("I136").Value = 2137.3823
$ws.Range("J136").Value = 17548
$ws.Range("K136").Value = 6412.146900000001
$ws.Range("L136").Value = 52644
$ws.Range("M136").Value = -3862.146900000001
$ws.Range("N136").Value = -57744

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1146.5883
$ws.Range("I122").Value = 1137
$ws.Range("J122").Value = 1300
$ws.Range("K122").Value = 3411
$ws.Range("L122").Value = 3900
$ws.Range("M122").Value = -961
$ws.Range("N122").Value = -8800
